# feat: add 2022-Q1 data
#
# Starting layout: [ "2021-Q3" (fund holdings), "总计" (totals/summary) ]
# Target layout:   [ "2021-Q3" (unchanged), "2022-Q1" (new fund holdings),
#                     "总计" (unchanged headers, one new row inserted on top) ]

$wb = $excel.ActiveWorkbook
$wsQ3 = $wb.Worksheets.Item(1)
$wsTotals = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by copying the "总计" sheet's format
#    (same header/border style) and inserting it right after "2021-Q3".
#    NOTE: once the copy is inserted, "总计" is pushed one slot further
#    down the tab order, so it must be re-fetched by name afterwards
#    instead of reusing the (now stale) $wsTotals handle.
# ---------------------------------------------------------------------
$wsTotals.Copy($null, $wsQ3)
$wsNew = $wb.Worksheets.Item(2)
$wsNew.Name = "2022-Q1"
$wsTotals = $wb.Worksheets.Item("总计")

# Extend the bold/bordered header style already on D1 across the new
# E1:H1 header cells before we overwrite their text.
$wsNew.Range("D1").Copy()
$wsNew.Range("E1:H1").PasteSpecial(-4122)

$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Row 2 holds the fund data. A2 (the rank/order number) keeps the style
# copied from the old sheet. The numeric-looking values are stored as
# plain text (matching the source data), so a leading apostrophe forces
# text entry and ClearFormats drops the leftover "quote prefix" style bit.
$wsNew.Range("B2").Value = "'968013"
$wsNew.Range("C2").Value = "施罗德亚洲高息股债基金M"
$wsNew.Range("D2").Value = "'297.64"
$wsNew.Range("E2").Value = "'57.54"
$wsNew.Range("F2").Value = "'1.29"
$wsNew.Range("G2").Value = "'3.8396"
$wsNew.Range("H2").Value = 7
$wsNew.Range("B2:G2").ClearFormats()

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new row above the existing
#    "2021-Q3" total row and fill it in with the 2022-Q1 totals.
# ---------------------------------------------------------------------
$wsTotals.Rows.Item(2).Insert()

# Give the new A2 the same style as A3 (the pre-existing ordinal cell).
$wsTotals.Range("A3").Copy()
$wsTotals.Range("A2").PasteSpecial(-4122)
# The row-insert leaves an inherited border style on B2:D2; drop it so
# these cells stay unstyled like the original data cells.
$wsTotals.Range("B2:D2").ClearFormats()

$wsTotals.Range("A2").Value = 0
$wsTotals.Range("B2").Value = "2022-Q1"
$wsTotals.Range("C2").Value = 1
$wsTotals.Range("D2").Value = 3.84

# Renumber the old row (now shifted to row 3).
$wsTotals.Range("A3").Value = 1
